# Generate Report for Handback
#
# The bd9676bb-e471-4eaf-be59-5f3047e621af.md file has now been handed
# back (for both zh-cn and de-de). Update the status / datetime columns
# and record the new "Latest Target File" / "Latest Handback File"
# hyperlinks + datetime on the per-language sheets, and roll the
# Overview sheet's status + zh-cn/de-de columns forward to reflect the
# handback.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 3 & 4 (bd9676bb-e471-4eaf-be59-5f3047e621af.md /
# f86d5912-d005-4275-bf77-855c5cbdcabf.md) move from "Ready for
# handoff" to "Handed back: in sync with en-US" for both zh-cn (B) and
# de-de (C). The "Latest Handoff Date" column (D) is unchanged.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack
$wsOverview.Range("B4").Value = $statusHandedBack
$wsOverview.Range("C4").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 3 & 4 for bd9676bb-e471-4eaf-be59-5f3047e621af.md
# Status -> Handed back; fill in Latest Target File (F) and Latest
# Handback File (G) hyperlinks; update Latest Handback DateTime (H).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("C4").Value = $statusHandedBack

$zhCnTargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d222b54163c23ce4f5e779a04d6c74c1d12fed64/e2e/bd9676bb-e471-4eaf-be59-5f3047e621af.md"
$zhCnHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d222b54163c23ce4f5e779a04d6c74c1d12fed64/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhCnTargetUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhCnHandbackUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), $zhCnTargetUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G4"), $zhCnHandbackUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.zh-cn.xlf")

$wsZhCn.Range("H3").Value = "2016-03-20 04:50:12"
$wsZhCn.Range("H4").Value = "2016-03-20 04:50:12"

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, different URLs/datetime.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("C4").Value = $statusHandedBack

$deDeTargetUrl = "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/515926459a9d86a9a0bd4621b06f597b596e831e/e2e/bd9676bb-e471-4eaf-be59-5f3047e621af.md"
$deDeHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/515926459a9d86a9a0bd4621b06f597b596e831e/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deDeTargetUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deDeHandbackUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), $deDeTargetUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G4"), $deDeHandbackUrl, "", "", "bd9676bb-e471-4eaf-be59-5f3047e621af.4a011a765f1d1584affc5d7dcd7dad8ecd2a4d34.de-de.xlf")

$wsDeDe.Range("H3").Value = "2016-03-20 04:50:27"
$wsDeDe.Range("H4").Value = "2016-03-20 04:50:27"

Write-Output "Handback report generated."
